$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Developer hourly rate (C6): 100 -> 192. Dependent formulas (D6, F6, B1)
# will recalculate automatically.
$ws.Range("C6").Value = 192

# Move the active selection to D13, matching the saved cursor position.
$ws.Range("D13").Select()
